$wb = $excel.ActiveWorkbook
$wsParam = $wb.Worksheets.Item("Параметры")

# --- Add the new company row to "Параметры" sheet ---
$wsParam.Range("A9").Value = "ПАО ""Газпром"""
$wsParam.Range("B9").Value = 7736050003
$wsParam.Range("C9").Value = 42858
$wsParam.Range("D9").Value = 43081

# Copy the date formatting from the row above onto the new date cells
$wsParam.Range("C7:D7").Copy()
$wsParam.Range("C9:D9").PasteSpecial(-4122)

# --- Grow the "inputTable" ListObject to cover the new row ---
$tbl = $wsParam.ListObjects.Item("inputTable")
$tbl.Resize($wsParam.Range("A4:D9"))

# --- Switch the active tab / selection to the "Параметры" sheet ---
$wsParam.Activate()
$wsParam.Range("B17").Select() | Out-Null
